# Generate Report for Archive
#
# The localization status for this item moved from "Ready for handoff" to
# "In Translation". That status string is shown in:
#   - Overview sheet: columns E ("zh-cn") and F ("de-de"), row 2
#   - zh-cn sheet:   column C ("Status"), row 2
#   - de-de sheet:   column C ("Status"), row 2
#
# Shortening the status text also shrinks the auto-fitted "Status" column,
# so the corresponding column widths are narrowed to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status text wherever it appears.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Narrow the "Status" columns to reflect the shorter text (auto-fit result).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
